# The candidate list on "tresquintos.cl" contained a duplicate entry for
# Araucanía / Eduardo Vicencio (row 69): it was accidentally listed twice
# (once with list "DC" and once, one row below, with list "PPD"). Remove
# the stray duplicate row so the remaining rows shift up one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate row; everything below moves up and the used range
# shrinks from A1:H88 to A1:H87.
$ws.Rows(69).Delete()

# Refresh the sheet's remembered sort range/conditions (previously
# A2:H89 / C / F) to match the new, one-row-shorter data extent.
$sort = $ws.Sort
$sort.SortFields.Add($ws.Range("C2:C88"))
$sort.SortFields.Add($ws.Range("F2:F88"))
$sort.SetRange($ws.Range("A2:H88"))
$sort.Apply()

# Leave the selection where the author ended up after the edit.
$ws.Range("F62:F69").Select()
